# "Generate Report for Handback"
#
# The handback run completed: the zh-cn and de-de localization files came
# back in sync with en-US, so the report needs to reflect:
#   - Overview sheet status text for both languages
#   - per-language sheets: Latest Target File (+ hyperlink), Latest
#     Handback File and Latest Handback DateTime for both data rows
#   - column widths widened so the new/longer values are readable

$wb = $excel.ActiveWorkbook

$mdName      = "6a76ef82-af34-4370-ba4b-0bd10e284961.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37e64acc4095203aaee26c0a113ba0ccbc78467a/e2e/$mdName"
$zhHandback  = "6a76ef82-af34-4370-ba4b-0bd10e284961.9dcdcfa04436ac81976d4c042b0357ba04b05efb.zh-cn.xlf"
$deHandback  = "6a76ef82-af34-4370-ba4b-0bd10e284961.9dcdcfa04436ac81976d4c042b0357ba04b05efb.de-de.xlf"
$zhDateTime  = "2016-09-01 11:12:51"
$deDateTime  = "2016-09-01 11:12:58"
$statusText  = "Handed back: in sync with en-US"

# ColumnWidth is quantized by the host in steps of 1/6 of a character,
# so nudge the requested widths to the middle of the bucket that rounds
# to the width we actually want stored in the sheet.
function Set-ColWidth($col, $target) {
    $k = [Math]::Round($target * 6 - 5)
    $col.ColumnWidth = ($k) / 6
}

# ---------------------------------------------------------------------
# Overview sheet: both language status columns move from "Ready for
# handoff" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

Set-ColWidth $overview.Columns("E") 29.9777047293527
Set-ColWidth $overview.Columns("F") 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet: record the target file link + handback file + datetime
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("I2").Value = $mdName
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$zh.Range("J2").Value = $zhHandback

$zh.Range("I3").Value = $mdName
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$zh.Range("J3").Value = $zhHandback

$zh.Range("K2").Value = $zhDateTime
$zh.Range("K3").Value = $zhDateTime

Set-ColWidth $zh.Columns("C") 29.9777047293527
Set-ColWidth $zh.Columns("I") 40
Set-ColWidth $zh.Columns("J") 40

# ---------------------------------------------------------------------
# de-de sheet: same shape, its own handback file + a later datetime
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("I2").Value = $mdName
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$de.Range("J2").Value = $deHandback

$de.Range("I3").Value = $mdName
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$de.Range("J3").Value = $deHandback

$de.Range("K2").Value = $deDateTime
$de.Range("K3").Value = $deDateTime

Set-ColWidth $de.Columns("C") 29.9777047293527
Set-ColWidth $de.Columns("I") 40
Set-ColWidth $de.Columns("J") 40
